# Doing Updates for Financials
# A new fiscal year (period ending 2018-12-31) is added as the first data
# column on the AFG financial-statements sheet. The existing data in
# columns D:J (covering the prior 7 periods) is shifted one column to the
# right (into E:K), a new column D is populated with the freshly reported
# figures, and a new (empty) trailing column L is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before D, pushing D:K -> E:L --------------
$ws.Columns("D:D").Insert()

# Copy number formats/styles from the (now shifted) column E into the
# newly inserted column D, restricted to the data range so we do not
# blow out the sheet dimension.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# --- 2. Populate the new column D with the latest reported period -----

# Income Statement ------------------------------------------------------
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 7169000
$ws.Range("D9").Value = 5885000
$ws.Range("D10").Value = 1284000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 19000
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 6468000
$ws.Range("D18").Value = 701000
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 911000
$ws.Range("D22").Value = 62000
$ws.Range("D23").Value = 639000
$ws.Range("D24").Value = 122000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 517000
$ws.Range("D27").Value = 530000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 530000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 530000

# Balance Sheet -----------------------------------------------------------
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 1515000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 5673000
$ws.Range("D44").Value = 0

# Row 45 ("Other Current Assets") is special-cased: the whole D:J run,
# which used to be the text placeholder "NA", is replaced with resolved
# numeric zeros, while the old J value (409000) shifts into K.
$ws.Range("D45:J45").Value = 0
$ws.Range("K45").Value = 409000

$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 51390000
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 246000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 63456000
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 3904000
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 1302000
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 58486000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 3588000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 4970000
$ws.Range("D77").Value = 0

# Cash Flow Statement ------------------------------------------------------
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 530000
$ws.Range("D83").Value = 210000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2083000
$ws.Range("D91").Value = -80000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -5350000
$ws.Range("D96").Value = -394000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 2444000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -823000

# --- 3. Housekeeping rows that have no data columns at all (pure section
#        headers / blank separators) only gained the new blank column as a
#        by-product of the column insert above and need no extra values.

$ws.Range("A5").Select()
